# Test for parsing NaN added: populate a new column G on Sheet3 with a
# "not available" header and two "NaN" values, matching the new lloq/value
# columns already present in columns C:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("G1").Value = "not available"
$ws.Range("G2").Value = "NaN"
$ws.Range("G3").Value = "NaN"

# Widen the new column similarly to how the author resized it by hand.
$ws.Columns.Item(7).ColumnWidth = 15.8

# Leave the selection where the author ended up after the edit.
$ws.Range("H12").Select()
